$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update conversion text in cell A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.39 = 12935.59 pesos`n✅ 12935.59 pesos = 3.38 = 962.61 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 295
$ws2.Range("O10").Value = 3816
$ws2.Range("N12").Value = 3829.99
$ws2.Range("O12").Value = 285.011
